$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A86").Value = "vertfieldgrad"
$ws.Range("B86").Value = "V"

$ws.Range("A87").Value = "LatticeScope_Ch3_Mean"
$ws.Range("B87").Style = "Normal"
